$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 18 from 45207 (2023-10-08)
# to 45208 (2023-10-09), matching the author's automatic date-refresh update.
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
